$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Talla" (size) formulas in column D, rows 2-8.
# Formula extracts the size code found after the "-" and the 2-digit color
# code, e.g. "100's:100-65XL" -> "XL"
# D2 is entered on its own, D3:D8 as one range so they share the formula,
# mirroring the layout already used by columns B and C.
$ws.Range("D2").Formula = '=MID(A2,FIND("-",A2,1)+3,3)'
$ws.Range("D3:D8").Formula = '=MID(A3,FIND("-",A3,1)+3,3)'

# Move the active selection to C11, matching the saved sheet view.
$ws.Range("C11").Select()
